# Apply updates described by the commit diff:
# 1) Update several odds values in rows 2-6 (existing matches, values re-priced).
# 2) Remove row 7 (Amazonas vs America MG) entirely; subsequent rows shift up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Independiente vs Union de Santa Fe) ---
$ws.Range("I2").Value  = 4.1
$ws.Range("J2").Value  = 3.2
$ws.Range("K2").Value  = 1.77
$ws.Range("N2").Value  = 5
$ws.Range("R2").Value  = 1.27
$ws.Range("X2").Value  = 8.5
$ws.Range("AC2").Value = 4.5
$ws.Range("AG2").Value = 7.5
$ws.Range("AP2").Value = 41
$ws.Range("AZ2").Value = 101

# --- Row 3 (Newells Old Boys vs Huracan) ---
$ws.Range("K3").Value = 1.69
$ws.Range("R3").Value = 1.22

# --- Row 4 (Bahia vs Sao Paulo) ---
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.53

# --- Row 5 (Botafogo RJ vs Vasco) ---
$ws.Range("Q5").Value = 1.73
$ws.Range("R5").Value = 2.08

# --- Row 6 (Internacional vs Criciuma) ---
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 13

# --- Remove row 7 (Amazonas vs America MG) so all following rows shift up ---
$ws.Rows.Item(7).Delete()
